$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new Student's Question / By / Date row data (row 3)
$ws.Range("B3").Value = "IS there any Sample projects like this?"
$ws.Range("C3").Value = "Ashkan"

# Match D3's date formatting to D2's (copy format only, then set the date value)
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("D3").Value = (Get-Date -Year 2022 -Month 9 -Day 21 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("D4").Select()
